# Generate Report for Handoff
# Update status "In Translation" -> "Ready for handoff" and refresh the
# handoff timestamps on the Overview, zh-cn, and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn status, de-de status, and latest handoff date
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-25-14 08:25:10"

# zh-cn sheet: status + latest handoff datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-14 08:25:07"

# de-de sheet: status + latest handoff datetime
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-14 08:25:10"
